# Sub-EMU "jitter" correction pass over the China-organization-location
# choropleth group (matches a re-render of the same map where a handful
# of coordinates drifted by a single EMU during rounding).
#
# PowerPoint's Shape.Left/Top/Width/Height are Single-precision (32-bit
# float) "points" properties; the literals below are the exact float32
# values whose point->EMU conversion lands on the target EMU figure, so
# that assigning them reproduces the precise before/after EMU deltas
# from the authoring diff without disturbing any sibling geometry.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

function Get-Item-ById($groupShape, $id) {
    return $groupShape.GroupItems.Item($id - 2)
}

# rc3 (id 3): ext cx 6400800 -> 6400799
$sh = Get-Item-ById $grp 3
$sh.Width = 503.99993896484375

# tx6 (id 6): off y 2463680 -> 2463679
$sh = Get-Item-ById $grp 6
$sh.Top = 193.990478515625

# pl11 (id 11): off y 2930827 -> 2930826
$sh = Get-Item-ById $grp 11
$sh.Top = 230.77371215820312

# pg13 (id 13): off x 4526315 -> 4526314
$sh = Get-Item-ById $grp 13
$sh.Left = 356.4026794433594

# pl14 (id 14): off x 4535231 -> 4535230
$sh = Get-Item-ById $grp 14
$sh.Left = 357.104736328125

# tx15 (id 15): off x 4089578 -> 4089577
$sh = Get-Item-ById $grp 15
$sh.Left = 322.0139465332031

# pg19 (id 19): off x 4461858 -> 4461857
$sh = Get-Item-ById $grp 19
$sh.Left = 351.32733154296875

# pg25 (id 25): off x 4562043 -> 4562042
$sh = Get-Item-ById $grp 25
$sh.Left = 359.2159118652344

# pg28 (id 28): off x 4712151 -> 4712150
$sh = Get-Item-ById $grp 28
$sh.Left = 371.03546142578125

# pl41 (id 41): off x 6537607 -> 6537608
$sh = Get-Item-ById $grp 41
$sh.Left = 514.7723388671875

# pg43 (id 43): off y 2619907 -> 2619906; ext cx 1239763 -> 1239764; ext cy 1070077 -> 1070078
$sh = Get-Item-ById $grp 43
$sh.Top = 206.2918243408203
$sh.Width = 97.61921691894531
$sh.Height = 84.25811767578125

Write-Host "done"
